$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, shifting existing rows 11:89 down to 12:90.
$ws.Rows("11:11").Insert()

# Populate the newly inserted row 11 with the new weekly data entry.
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44490
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112040
$ws.Cells.Item(11, 7).Value = "Cilantro"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 3200
$ws.Cells.Item(11, 11).Value = 1300
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 13).Value = 1400
$ws.Cells.Item(11, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 933
$ws.Cells.Item(11, 17).Value = 1.5
$ws.Cells.Item(11, 18).Value = "Hortaliza"
